$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3535
$ws.Range("J19").Value = 6388
$ws.Range("L19").Value = 6388
$ws.Range("N19").Value = -6738
$ws.Range("H39").Value = 394.875
$ws.Range("I39").Value = 149
$ws.Range("J39").Value = 935.8
$ws.Range("K39").Value = 447
$ws.Range("L39").Value = 2807.4
$ws.Range("M39").Value = -151
$ws.Range("N39").Value = -3399.4
$ws.Range("H48").Value = 4600
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = -2708
$ws.Range("N48").Value = -30584
$ws.Range("H51").Value = 8628.571
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8628.571
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 8628.571
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -9596.571
$ws.Range("H56").Value = 4600
$ws.Range("I56").Value = 1000
$ws.Range("J56").Value = 10000
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 30000
$ws.Range("M56").Value = -2466
$ws.Range("N56").Value = -31068
$ws.Range("H95").Value = 22370.8
$ws.Range("J95").Value = 22370.8
$ws.Range("L95").Value = 22370.8
$ws.Range("N95").Value = -27862.8
$ws.Range("H98").Value = 5785.5557
$ws.Range("I98").Value = 5785.5557
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 5785.5557
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -4287.5557
$ws.Range("N98").ClearContents()
$ws.Range("H112").Value = 1508.7727
$ws.Range("J112").Value = 1508.7727
$ws.Range("L112").Value = 4526.3181
$ws.Range("N112").Value = -6742.3181
$ws.Range("H122").Value = 5785.5557
$ws.Range("I122").Value = 5785.5557
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17356.6671
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14906.6671
$ws.Range("N122").ClearContents()
$ws.Range("H137").Value = 2745.2415
$ws.Range("J137").Value = 2807.9443
$ws.Range("L137").Value = 8423.832900000001
$ws.Range("N137").Value = -13523.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 8380.571
$ws.Range("I4").Value = 209.66667
$ws.Range("K4").Value = 209.66667
$ws.Range("M4").Value = -93.66667000000001
$ws.Range("H14").Value = 401.125
$ws.Range("I14").Value = 251.5
$ws.Range("K14").Value = 251.5
$ws.Range("M14").Value = -76.5
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -771
$ws.Range("N19").ClearContents()
$ws.Range("H29").Value = 2002468.8
$ws.Range("I29").Value = 2502961
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 2502961
$ws.Range("L29").Value = 500
$ws.Range("M29").Value = -2502653
$ws.Range("N29").Value = -1116
$ws.Range("H30").Value = 1345.6
$ws.Range("I30").Value = 1886
$ws.Range("J30").Value = 535
$ws.Range("K30").Value = 1886
$ws.Range("L30").Value = 535
$ws.Range("M30").Value = -1736
$ws.Range("N30").Value = -835
$ws.Range("H35").Value = 2913.25
$ws.Range("I35").Value = 2913.25
$ws.Range("K35").Value = 2913.25
$ws.Range("M35").Value = -2507.25
$ws.Range("H53").Value = 8000
$ws.Range("I53").Value = 8000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 8000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -7318
$ws.Range("N53").ClearContents()
$ws.Range("H122").Value = 1837.85
$ws.Range("I122").Value = 1583.8667
$ws.Range("K122").Value = 4751.6001
$ws.Range("M122").Value = -2301.6001
$ws.Range("H132").Value = 2496.2856
$ws.Range("I132").Value = 1420.75
$ws.Range("K132").Value = 4262.25
$ws.Range("M132").Value = -1732.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2725.6
$ws.Range("I36").Value = 2445.6667
$ws.Range("J36").Value = 3145.5
$ws.Range("K36").Value = 2445.6667
$ws.Range("L36").Value = 3145.5
$ws.Range("M36").Value = -1911.6667
$ws.Range("N36").Value = -4213.5
$ws.Range("H39").Value = 3011
$ws.Range("J39").Value = 3222
$ws.Range("L39").Value = 3222
$ws.Range("N39").Value = -4000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3006.2532
$ws.Range("I31").Value = 1922.8
$ws.Range("J31").Value = 3507.8518
$ws.Range("K31").Value = 1922.8
$ws.Range("L31").Value = 3507.8518
$ws.Range("M31").Value = -1627.8
$ws.Range("N31").Value = -4097.8518
$ws.Range("H34").Value = 3006.2532
$ws.Range("I34").Value = 1922.8
$ws.Range("J34").Value = 3507.8518
$ws.Range("K34").Value = 1922.8
$ws.Range("L34").Value = 3507.8518
$ws.Range("M34").Value = -1720.8
$ws.Range("N34").Value = -3911.8518
$ws.Range("H122").Value = 2406.8572
$ws.Range("I122").Value = 1780
$ws.Range("J122").Value = 3974
$ws.Range("K122").Value = 5340
$ws.Range("L122").Value = 11922
$ws.Range("M122").Value = -2890
$ws.Range("N122").Value = -16822

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1952.7037
$ws.Range("J23").Value = 3008.375
$ws.Range("L23").Value = 9025.125
$ws.Range("N23").Value = -9495.125
$ws.Range("H40").Value = 1611
$ws.Range("H132").Value = 2738.4285
$ws.Range("I132").Value = 1312.1666
$ws.Range("J132").Value = 2872.1406
$ws.Range("K132").Value = 11809.4994
$ws.Range("L132").Value = 25849.2654
$ws.Range("M132").Value = -9279.499400000001
$ws.Range("N132").Value = -30909.2654

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 991.2222
$ws.Range("I97").Value = 529
$ws.Range("J97").Value = 1569
$ws.Range("K97").Value = 529
$ws.Range("L97").Value = 1569
$ws.Range("M97").Value = -33
$ws.Range("N97").Value = -2561
$ws.Range("H122").Value = 2364.45
$ws.Range("I122").Value = 1814.5385
$ws.Range("K122").Value = 5443.6155
$ws.Range("M122").Value = -2993.6155
$ws.Range("H134").Value = 18331.5
$ws.Range("J134").Value = 18331.5
$ws.Range("L134").Value = 54994.5
$ws.Range("N134").Value = -60064.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2711.8845
$ws.Range("I46").Value = 2482.3333
$ws.Range("J46").Value = 2780.75
$ws.Range("K46").Value = 2482.3333
$ws.Range("L46").Value = 2780.75
$ws.Range("M46").Value = -2294.3333
$ws.Range("N46").Value = -3156.75
$ws.Range("H93").Value = 11907836
$ws.Range("I93").Value = 13891892
$ws.Range("J93").Value = 3501.75
$ws.Range("K93").Value = 13891892
$ws.Range("L93").Value = 3501.75
$ws.Range("M93").Value = -13890644
$ws.Range("N93").Value = -5997.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8400
$ws.Range("I62").Value = 8000
$ws.Range("K62").Value = 8000
$ws.Range("M62").Value = -7376
$ws.Range("H65").Value = 8400
$ws.Range("I65").Value = 8000
$ws.Range("K65").Value = 40000
$ws.Range("M65").Value = -36880
$ws.Range("H81").Value = 3415.8333
$ws.Range("I81").Value = 2291.6428
$ws.Range("J81").Value = 4989.7
$ws.Range("K81").Value = 4583.2856
$ws.Range("L81").Value = 9979.4
$ws.Range("M81").Value = -3522.2856
$ws.Range("N81").Value = -12101.4
$ws.Range("H84").Value = 3415.8333
$ws.Range("I84").Value = 2291.6428
$ws.Range("J84").Value = 4989.7
$ws.Range("K84").Value = 22916.428
$ws.Range("L84").Value = 49897
$ws.Range("M84").Value = -17612.428
$ws.Range("N84").Value = -60505
$ws.Range("H96").Value = 3741.6667
$ws.Range("I96").Value = 3112.5
$ws.Range("K96").Value = 3112.5
$ws.Range("M96").Value = -1739.5
$ws.Range("H113").Value = 6991.7896
$ws.Range("I113").Value = 7016.2666
$ws.Range("J113").Value = 6900
$ws.Range("K113").Value = 21048.7998
$ws.Range("L113").Value = 20700
$ws.Range("M113").Value = -18878.7998
$ws.Range("N113").Value = -25040
$ws.Range("H122").Value = 1738.9333
$ws.Range("I122").Value = 1511.2084
$ws.Range("K122").Value = 4533.6252
$ws.Range("M122").Value = -2083.6252
$ws.Range("H136").Value = 3633.077
$ws.Range("I136").Value = 3654.9565
$ws.Range("K136").Value = 10964.8695
$ws.Range("M136").Value = -8414.869499999999
